$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6470
$ws.Range("L3").Value = 6978
$ws.Range("L4").Value = 1740
$ws.Range("L5").Value = 412
$ws.Range("L6").Value = 5736
$ws.Range("L7").Value = 21336

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 497
$ws.Range("L4").Value = 98
$ws.Range("L7").Value = 1410

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 339
$ws.Range("L7").Value = 959

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 249
$ws.Range("L6").Value = 210
$ws.Range("L7").Value = 819

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 156
$ws.Range("L3").Value = 133
$ws.Range("L7").Value = 421

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 77
$ws.Range("L7").Value = 677
$ws.Range("L8").Value = 1410
$ws.Range("L9").Value = 122
$ws.Range("L12").Value = 50
$ws.Range("L15").Value = 182
$ws.Range("L19").Value = 588
$ws.Range("L23").Value = 224
$ws.Range("L29").Value = 1193
$ws.Range("L33").Value = 959
$ws.Range("L37").Value = 819
$ws.Range("L42").Value = 674
$ws.Range("L43").Value = 159
$ws.Range("L52").Value = 453
$ws.Range("L54").Value = 459
$ws.Range("L55").Value = 227
$ws.Range("L57").Value = 72
$ws.Range("L60").Value = 145
$ws.Range("L65").Value = 421
$ws.Range("L67").Value = 743
$ws.Range("L73").Value = 169
$ws.Range("L76").Value = 337
$ws.Range("L79").Value = 591
$ws.Range("L84").Value = 205
$ws.Range("L85").Value = 1057
$ws.Range("L88").Value = 226
$ws.Range("L89").Value = 288
$ws.Range("L90").Value = 227
$ws.Range("L93").Value = 108
$ws.Range("L94").Value = 259
$ws.Range("L98").Value = 112
$ws.Range("L101").Value = 21336

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 210
$ws.Range("L3").Value = 288
$ws.Range("L4").Value = 50
$ws.Range("L6").Value = 174
$ws.Range("L7").Value = 743

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 67
$ws.Range("L5").Value = 5
$ws.Range("L7").Value = 205

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L5").Value = 5
$ws.Range("L7").Value = 459

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 364
$ws.Range("L4").Value = 65
$ws.Range("L7").Value = 1193

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 210
$ws.Range("L3").Value = 179
$ws.Range("L6").Value = 162
$ws.Range("L7").Value = 588

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L4").Value = 44
$ws.Range("L6").Value = 152
$ws.Range("L7").Value = 337

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 232
$ws.Range("L6").Value = 191
$ws.Range("L7").Value = 674

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 66
$ws.Range("L7").Value = 227

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 183
$ws.Range("L3").Value = 192
$ws.Range("L7").Value = 591

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 38
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 217
$ws.Range("L6").Value = 162
$ws.Range("L7").Value = 677

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 259

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 67
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 75
$ws.Range("L7").Value = 288

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 77
$ws.Range("L7").Value = 227

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 159

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 316
$ws.Range("L3").Value = 437
$ws.Range("L6").Value = 222
$ws.Range("L7").Value = 1057

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 141
$ws.Range("L7").Value = 453

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 50
